# test #8: Modelos finales de tests tras el filtro de celdas vacias de los
# excel importados.
#
# The imported census sample (testImport3.xlsx) gets its "empty cell"
# placeholders replaced with real numeric data, and a couple of existing
# numbers are swapped out for their final values. The active selection on
# the sheet also moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: voting_id 7 -> 87885, voter_id 363 -> 3863
$ws.Range("A2").Value = 87885
$ws.Range("D2").Value = 3863

# Row 3: voting_id was empty -> 52414, voter_id 3 -> 8683
$ws.Range("A3").Value = 52414
$ws.Range("D3").Value = 8683

# Row 4: voting_id 9 -> 2, voter_id was empty -> 1
$ws.Range("A4").Value = 2
$ws.Range("D4").Value = 1

# Row 5: voting_id 56 -> 68688
$ws.Range("A5").Value = 68688

# The sheet's saved active cell/selection moves from D2 to C7
$ws.Range("C7").Select()
